$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff represents a cyclic shift of the "Fecha" (D), "Volumen" (M),
# "Precio máximo" (O), "Precio promedio ponderado" (P) and "Precio $/Kg" (S)
# values across rows 2, 3 and 4: row2 <- old row3, row3 <- old row4,
# row4 <- old row2. "Precio mínimo" (N) stays unchanged (180000 for all rows).

$oldD2 = $ws.Range("D2").Value2
$oldD3 = $ws.Range("D3").Value2
$oldD4 = $ws.Range("D4").Value2

$oldM2 = $ws.Range("M2").Value2
$oldM3 = $ws.Range("M3").Value2
$oldM4 = $ws.Range("M4").Value2

$oldO2 = $ws.Range("O2").Value2
$oldO3 = $ws.Range("O3").Value2
$oldO4 = $ws.Range("O4").Value2

$oldP2 = $ws.Range("P2").Value2
$oldP3 = $ws.Range("P3").Value2
$oldP4 = $ws.Range("P4").Value2

$oldS2 = $ws.Range("S2").Value2
$oldS3 = $ws.Range("S3").Value2
$oldS4 = $ws.Range("S4").Value2

$ws.Range("D2").Value2 = $oldD3
$ws.Range("D3").Value2 = $oldD4
$ws.Range("D4").Value2 = $oldD2

$ws.Range("M2").Value2 = $oldM3
$ws.Range("M3").Value2 = $oldM4
$ws.Range("M4").Value2 = $oldM2

$ws.Range("O2").Value2 = $oldO3
$ws.Range("O3").Value2 = $oldO4
$ws.Range("O4").Value2 = $oldO2

$ws.Range("P2").Value2 = $oldP3
$ws.Range("P3").Value2 = $oldP4
$ws.Range("P4").Value2 = $oldP2

$ws.Range("S2").Value2 = $oldS3
$ws.Range("S3").Value2 = $oldS4
$ws.Range("S4").Value2 = $oldS2
